# Daily auto-update of the price-data sheet: prepend a new row for the
# latest date (2025-12-25), pushing every existing row down by one so the
# table keeps its newest-first ordering. The numeric columns carry over
# the same (unchanged) constant values used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 2..35 down to 3..36, opening up a blank row 2 for today's entry.
$ws.Rows.Item(2).Insert()

# Force column A to be treated as plain text so the date string isn't
# auto-converted into a date serial number.
$ws.Cells.Item(2, 1).NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "2025-12-25"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610

# Drop the formatting picked up from the row-insert (and the NumberFormat
# tweak above) so the new row matches the unformatted look of every other
# data row in the table.
$ws.Range("A2:D2").ClearFormats()
